$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewLoanInput")
$ws.Range("B2").Value = "3500-RBI-EI-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-1-ONTIME-DISBURSE-FEE-%APR-AMT-Reg-PERIODIC"
$ws.Range("B8").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B2").Select() | Out-Null
